# "fix upload excel admin" - reset the bad/test row (row 2) that was
# uploaded incorrectly: mark it as "Tidak Berpenghuni" and clear out all of
# the questionnaire answer columns (which were populated with stray
# leftover values) back to 0, keeping only the identifying columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: status changes from "Berpenghuni" to "Tidak Berpenghuni"
$ws.Range("A2").Value = "Tidak Berpenghuni"

# B2, C2, D2, L2, M2, N2, BG2, BH2 stay as-is (unchanged)

# E2..K2 reset to 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# O2..K2 (O..BE, skipping L/M/N which are kept) reset to 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = 0
$ws.Range("AP2").Value = 0
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = 0
$ws.Range("AS2").Value = 0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0
$ws.Range("AW2").Value = 0
$ws.Range("AX2").Value = 0
$ws.Range("AY2").Value = 0
$ws.Range("AZ2").Value = 0
$ws.Range("BA2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = 0
$ws.Range("BE2").Value = 0

# BG2 / BH2 stay as-is (unchanged)

# Move the visible selection to H2 (also scrolls the view back to A1,
# dropping the old topLeftCell="BC1" scroll position)
$ws.Range("A1").Select()
$ws.Range("H2").Select()
